$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 327, shifting the existing rows 327:341 down to 328:342
$ws.Rows(327).Insert()

# Populate the newly inserted row 327 with the new weekly price record
$ws.Cells.Item(327, 1).Value  = 4
$ws.Cells.Item(327, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(327, 3).Value  = "Los Lagos"
$ws.Cells.Item(327, 4).Value  = 44931
$ws.Cells.Item(327, 5).Value  = 10
$ws.Cells.Item(327, 6).Value  = "Fruta"
$ws.Cells.Item(327, 7).Value  = 100101
$ws.Cells.Item(327, 8).Value  = "Berries"
$ws.Cells.Item(327, 9).Value  = 100112025
$ws.Cells.Item(327, 10).Value = "Frutilla"
$ws.Cells.Item(327, 11).Value = "Sin especificar"
$ws.Cells.Item(327, 12).Value = "Primera"
$ws.Cells.Item(327, 13).Value = 800
$ws.Cells.Item(327, 14).Value = 8000
$ws.Cells.Item(327, 15).Value = 9000
$ws.Cells.Item(327, 16).Value = 8500
$ws.Cells.Item(327, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(327, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(327, 19).Value = 1214
$ws.Cells.Item(327, 20).Value = 7

# Match the date cell number format used by the other rows in column D
$ws.Cells.Item(327, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
